$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C24").Value = "17;20"

$ws.Application.ActiveWindow.Zoom = 137
$ws.Range("D27").Select()
